$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain Text so the stored cell type matches the source data.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.156.94"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "1.827.54"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +3.68%  "
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "43.01"
$ws.Range("E8").Value = "  +6.15%  "
$ws.Range("E9").Value = "  +6.82%  "
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "2.097.08"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "1.830.84"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "11.20"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").Value = "0.666"
$ws.Range("E15").Value = "  +5.68%  "
$ws.Range("D16").Value = "4.70"
$ws.Range("E16").Value = "  +6.86%  "
$ws.Range("D17").Value = "35.148.18"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "70.11"
$ws.Range("E18").Value = "  +3.92%  "
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").Value = "240.04"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "11.85"
$ws.Range("E21").Value = "  +7.04%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.60"
$ws.Range("E22").Value = "  +12.72%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.01"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("D25").Value = "171.55"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").Value = "17.59"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  +30.21%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "3.344.32"
$ws.Range("E31").Value = "  +37.64%  "
$ws.Range("D32").Value = "0.0555"
$ws.Range("E32").Value = "  +8.13%  "
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").Value = "4.01"
$ws.Range("E34").Value = "  +4.69%  "
$ws.Range("D35").Value = "1.80"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").Value = "93.87"
$ws.Range("E36").Value = "  +11.53%  "
$ws.Range("D37").Value = "0.682"
$ws.Range("E37").Value = "  +5.77%  "
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("D39").Value = "1.326.06"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "15.00"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.27"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").Value = "  +5.96%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("E47").Value = "  +7.84%  "
$ws.Range("D48").Value = "0.0509"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "2.007.50"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "100.86"
$ws.Range("E51").Value = "  +0.12%  "

# Restore default styling on the forced-text cells (Excel created a new
# number-format style when we set NumberFormat; drop it so cell styling
# matches the original "no explicit style" state).
$ws.Range("D8").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
